$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet - the "Ready for handoff" status text is shared with the
# per-file sheets below, so it flips to the new status here too.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (2e932acd-...)
$ws.Range("C2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md", "", "", "2e932acd-e47f-4f3c-8372-e61745a5bd03.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94b2af89cfdad9fae4b802d0f17b1861ba7d8d6c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf", "", "", "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf")
$ws.Range("H2").Value = "2016-03-25 10:53:43"

# Row 3 (77232830-...)
$ws.Range("C3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/77232830-5d71-4781-bf09-c9d381f516af.md", "", "", "77232830-5d71-4781-bf09-c9d381f516af.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94b2af89cfdad9fae4b802d0f17b1861ba7d8d6c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf", "", "", "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf")
$ws.Range("H3").Value = "2016-03-25 10:53:43"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (2e932acd-...)
$ws.Range("C2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md", "", "", "2e932acd-e47f-4f3c-8372-e61745a5bd03.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a0097503ccc873554da1958355484159f060a44/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf", "", "", "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf")
$ws.Range("H2").Value = "2016-03-25 10:53:58"

# Row 3 (77232830-...)
$ws.Range("C3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/77232830-5d71-4781-bf09-c9d381f516af.md", "", "", "77232830-5d71-4781-bf09-c9d381f516af.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a0097503ccc873554da1958355484159f060a44/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf", "", "", "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf")
$ws.Range("H3").Value = "2016-03-25 10:53:58"
